$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data row (STT/CSR/Certificate/Cert Chain sample row) so the
# stale shared strings are dropped and the sheet starts clean.
$ws.Range("A1:D2").ClearContents()

# Write the new header row (A1:H1) with the CSR export field names.
$ws.Range("A1").Value = "EmailAddress"
$ws.Range("B1").Value = "TelephoneNumber"
$ws.Range("C1").Value = "Locality"
$ws.Range("D1").Value = "StateProvince"
$ws.Range("E1").Value = "Country"
$ws.Range("F1").Value = "CustomerPhoneNumber"
$ws.Range("G1").Value = "CustomerEmail"
$ws.Range("H1").Value = "CSR"

# Resize the columns to match the new 8-column layout.
$ws.Columns.Item(1).ColumnWidth = 22.666666666666668
$ws.Columns.Item(2).ColumnWidth = 18.666666666666668
$ws.Columns.Item(3).ColumnWidth = 14.833333333333334
$ws.Columns.Item(4).ColumnWidth = 18.666666666666668
$ws.Columns.Item(5).ColumnWidth = 14.833333333333334
$ws.Columns.Item(6).ColumnWidth = 34.333333333333336
$ws.Columns.Item(7).ColumnWidth = 30.333333333333332
$ws.Columns.Item(8).ColumnWidth = 22.666666666666668

# Move the active selection to E21, matching the saved view state.
$ws.Range("E21").Select() | Out-Null
